$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 81, pushing existing rows 81-90 down to 82-91.
$ws.Rows.Item(81).Insert()

# Populate the new row 81 with the new record (same constant columns as
# the rest of the data set, new values for the variable columns).
$ws.Cells.Item(81, 1).Value = 7
$ws.Cells.Item(81, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(81, 3).Value = "Ñuble"
$ws.Cells.Item(81, 4).Value = 44782
$ws.Cells.Item(81, 5).Value = 16
$ws.Cells.Item(81, 6).Value = 100112021
$ws.Cells.Item(81, 7).Value = "Ají"
$ws.Cells.Item(81, 8).Value = "Inferno"
$ws.Cells.Item(81, 9).Value = "Primera"
$ws.Cells.Item(81, 10).Value = 60
$ws.Cells.Item(81, 11).Value = 15000
$ws.Cells.Item(81, 12).Value = 16000
$ws.Cells.Item(81, 13).Value = 15500
$ws.Cells.Item(81, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(81, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(81, 16).Value = 1033
$ws.Cells.Item(81, 17).Value = 15
$ws.Cells.Item(81, 18).Value = "Hortaliza"
